# smol refactoring code for generate_report
#
# Adds two new report sheets at the end of the workbook, cloned from the
# most recent existing report sheet (same layout/content/styles), then
# renamed to the new report timestamps.

$wb = $excel.ActiveWorkbook

# The template to clone: the last existing report sheet.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)

# First new sheet: report_2024_05_23-180202
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "report_2024_05_23-180202"

# Second new sheet: report_2024_05_23-180444
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "report_2024_05_23-180444"
